# Update the "Förändrad" date column (C2:C11) from 2023-09-11 (45180)
# to 2023-09-12 (45181), matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 11; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45180) {
        $cell.Value2 = 45181
    }
}
